# Scheduled runner update: refresh currentAveragePrice / Leve price / profit
# columns (H:N) for a handful of leve rows across the crafting-class sheets,
# reflecting newer market-board price data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 356.22726
$ws.Range("I9").Value = 338.84616
$ws.Range("K9").Value = 338.84616
$ws.Range("M9").Value = -169.84616
$ws.Range("H15").Value = 1534.3903
$ws.Range("I15").Value = 1534.3903
$ws.Range("K15").Value = 4603.1709
$ws.Range("M15").Value = -4434.1709
$ws.Range("H19").Value = 1137.5714
$ws.Range("I19").Value = 1211.5714
$ws.Range("J19").Value = 989.5714
$ws.Range("K19").Value = 1211.5714
$ws.Range("L19").Value = 989.5714
$ws.Range("M19").Value = -1036.5714
$ws.Range("N19").Value = -1339.5714
$ws.Range("H112").Value = 1391807.8
$ws.Range("I112").Value = 2156.3333
$ws.Range("J112").Value = 1590329.4
$ws.Range("K112").Value = 6468.999899999999
$ws.Range("L112").Value = 4770988.199999999
$ws.Range("M112").Value = -5360.999899999999
$ws.Range("N112").Value = -4773204.199999999
$ws.Range("H132").Value = 1703.238
$ws.Range("I132").Value = 1643.25
$ws.Range("K132").Value = 4929.75
$ws.Range("M132").Value = -2399.75
$ws.Range("H134").Value = 111361.25
$ws.Range("J134").Value = 111361.25
$ws.Range("L134").Value = 111361.25
$ws.Range("N134").Value = -121501.25
$ws.Range("H135").Value = 1800.7693
$ws.Range("I135").Value = 1004.1667
$ws.Range("K135").Value = 9037.5003
$ws.Range("M135").Value = -6502.5003
$ws.Range("H138").Value = 2487.33
$ws.Range("J138").Value = 2745.4578
$ws.Range("L138").Value = 8236.3734
$ws.Range("N138").Value = -18516.3734
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 4777.4443
$ws.Range("I45").Value = 4124.25
$ws.Range("J45").Value = 5300
$ws.Range("K45").Value = 4124.25
$ws.Range("L45").Value = 5300
$ws.Range("M45").Value = -3747.25
$ws.Range("N45").Value = -6054
$ws.Range("H76").Value = 16500
$ws.Range("I76").Value = 16500
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 16500
$ws.Range("L76").ClearContents()
$ws.Range("N76").Value = 0
$ws.Range("M76").Value = -16162
$ws.Range("H79").Value = 16500
$ws.Range("I79").Value = 16500
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 16500
$ws.Range("L79").ClearContents()
$ws.Range("N79").Value = 0
$ws.Range("M79").Value = -15330
$ws.Range("H102").Value = 1876.1538
$ws.Range("I102").Value = 1720
$ws.Range("K102").Value = 1720
$ws.Range("M102").Value = -98
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 810.8
$ws.Range("I64").Value = 250
$ws.Range("J64").Value = 951
$ws.Range("K64").Value = 250
$ws.Range("L64").Value = 951
$ws.Range("M64").Value = -25
$ws.Range("N64").Value = -1401
$ws.Range("H67").Value = 810.8
$ws.Range("I67").Value = 250
$ws.Range("J67").Value = 951
$ws.Range("K67").Value = 250
$ws.Range("L67").Value = 951
$ws.Range("M67").Value = 530
$ws.Range("N67").Value = -2511
$ws.Range("H105").Value = 2332.7778
$ws.Range("I105").Value = 1663.1818
$ws.Range("K105").Value = 1663.1818
$ws.Range("M105").Value = 83.81819999999993
$ws.Range("H132").Value = 101222.5
$ws.Range("J132").Value = 109963.336
$ws.Range("L132").Value = 109963.336
$ws.Range("N132").Value = -120083.336
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 295.45456
$ws.Range("I22").Value = 255.5
$ws.Range("K22").Value = 255.5
$ws.Range("M22").Value = 94.5
$ws.Range("H105").Value = 1468.9259
$ws.Range("I105").Value = 644.73334
$ws.Range("K105").Value = 644.73334
$ws.Range("M105").Value = 1102.26666
$ws.Range("H134").Value = 2500
$ws.Range("I134").Value = 2500
$ws.Range("K134").Value = 7500
$ws.Range("M134").Value = -4965
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 397.83334
$ws.Range("I2").Value = 375
$ws.Range("J2").Value = 409.25
$ws.Range("K2").Value = 2250
$ws.Range("L2").Value = 2455.5
$ws.Range("M2").Value = -2137
$ws.Range("N2").Value = -2681.5
$ws.Range("H107").Value = 986.875
$ws.Range("J107").Value = 836.1177
$ws.Range("L107").Value = 2508.3531
$ws.Range("N107").Value = -6348.3531
$ws.Range("H133").Value = 3410.5334
$ws.Range("I133").Value = 3166
$ws.Range("J133").Value = 5000
$ws.Range("K133").Value = 9498
$ws.Range("L133").Value = 15000
$ws.Range("M133").Value = -4438
$ws.Range("N133").Value = -25120
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2208
$ws.Range("I102").Value = 2028.8
$ws.Range("K102").Value = 2028.8
$ws.Range("M102").Value = -406.8
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4241.9585
$ws.Range("J40").Value = 5708.3335
$ws.Range("L40").Value = 5708.3335
$ws.Range("N40").Value = -5980.3335
$ws.Range("H42").Value = 16248.75
$ws.Range("J42").Value = 19985
$ws.Range("L42").Value = 19985
$ws.Range("N42").Value = -21111
$ws.Range("H48").Value = 26374
$ws.Range("I48").Value = 26374
$ws.Range("K48").Value = 26374
$ws.Range("M48").Value = -25713
$ws.Range("H49").Value = 16248.75
$ws.Range("J49").Value = 19985
$ws.Range("L49").Value = 19985
$ws.Range("N49").Value = -20279
$ws.Range("H50").Value = 30076
$ws.Range("I50").Value = 30076
$ws.Range("K50").Value = 30076
$ws.Range("M50").Value = -29439
$ws.Range("H55").Value = 380.44446
$ws.Range("I55").Value = 319.14285
$ws.Range("J55").Value = 595
$ws.Range("K55").Value = 319.14285
$ws.Range("L55").Value = 595
$ws.Range("M55").Value = -146.14285
$ws.Range("N55").Value = -941
$ws.Range("H56").Value = 18735.5
$ws.Range("I56").Value = 18735.5
$ws.Range("K56").Value = 18735.5
$ws.Range("M56").Value = -18044.5
$ws.Range("H82").Value = 2187.96
$ws.Range("I82").Value = 2627.5
$ws.Range("J82").Value = 1406.5555
$ws.Range("K82").Value = 2627.5
$ws.Range("L82").Value = 1406.5555
$ws.Range("M82").Value = -2266.5
$ws.Range("N82").Value = -2128.5555
$ws.Range("H85").Value = 2187.96
$ws.Range("I85").Value = 2627.5
$ws.Range("J85").Value = 1406.5555
$ws.Range("K85").Value = 2627.5
$ws.Range("L85").Value = 1406.5555
$ws.Range("M85").Value = -1379.5
$ws.Range("N85").Value = -3902.5555
$ws.Range("H93").Value = 2416.36
$ws.Range("I93").Value = 1973.2778
$ws.Range("J93").Value = 3555.7144
$ws.Range("K93").Value = 1973.2778
$ws.Range("L93").Value = 3555.7144
$ws.Range("M93").Value = -725.2778000000001
$ws.Range("N93").Value = -6051.7144
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2083.5518
$ws.Range("I132").Value = 1952.8
$ws.Range("K132").Value = 5858.4
$ws.Range("M132").Value = -3328.4
